$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 176; this shifts rows 176:222 down to 177:223
$ws.Rows("176:176").Insert()

# Populate the newly inserted row 176 with the new record's data
$ws.Range("A176").Value = 5
$ws.Range("B176").Value = "Macroferia Regional de Talca"
$ws.Range("C176").Value = "Maule"
$ws.Range("D176").Value = 44642
$ws.Range("E176").Value = 7
$ws.Range("F176").Value = 100112021
$ws.Range("G176").Value = "Ají"
$ws.Range("H176").Value = "Cristal"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 150
$ws.Range("K176").Value = 15000
$ws.Range("L176").Value = 15000
$ws.Range("M176").Value = 15000
$ws.Range("N176").Value = "`$/saco 25 kilos"
$ws.Range("O176").Value = "Región del Maule"
$ws.Range("P176").Value = 600
$ws.Range("Q176").Value = 25
$ws.Range("R176").Value = "Hortaliza"
